# Auto-generated edit script: applies scheduled market-data refresh values
# to the Leve profit tables (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 3594.75  # ALC!H103
$ws.Cells.Item(103, 9).Value = 4071.6  # ALC!I103
$ws.Cells.Item(103, 10).Value = 2800  # ALC!J103
$ws.Cells.Item(103, 11).Value = 12214.8  # ALC!K103
$ws.Cells.Item(103, 12).Value = 8400  # ALC!L103
$ws.Cells.Item(103, 13).Value = -11628.8  # ALC!M103
$ws.Cells.Item(103, 14).Value = -9572  # ALC!N103
$ws.Cells.Item(112, 8).Value = 1690.2  # ALC!H112
$ws.Cells.Item(112, 10).Value = 1728.0555  # ALC!J112
$ws.Cells.Item(112, 12).Value = 5184.166499999999  # ALC!L112
$ws.Cells.Item(112, 14).Value = -7400.166499999999  # ALC!N112
$ws.Cells.Item(113, 8).Value = 33337252  # ALC!H113
$ws.Cells.Item(113, 9).Value = 66668828  # ALC!I113
$ws.Cells.Item(113, 10).Value = 5678.2  # ALC!J113
$ws.Cells.Item(113, 11).Value = 66668828  # ALC!K113
$ws.Cells.Item(113, 12).Value = 5678.2  # ALC!L113
$ws.Cells.Item(113, 13).Value = -66665574  # ALC!M113
$ws.Cells.Item(113, 14).Value = -12186.2  # ALC!N113
$ws.Cells.Item(137, 8).Value = 1937.1765  # ALC!H137
$ws.Cells.Item(137, 9).Value = 1927.375  # ALC!I137
$ws.Cells.Item(137, 10).Value = 2094  # ALC!J137
$ws.Cells.Item(137, 11).Value = 5782.125  # ALC!K137
$ws.Cells.Item(137, 12).Value = 6282  # ALC!L137
$ws.Cells.Item(137, 13).Value = -3232.125  # ALC!M137
$ws.Cells.Item(137, 14).Value = -11382  # ALC!N137
$ws.Cells.Item(138, 8).Value = 3866.2222  # ALC!H138
$ws.Cells.Item(138, 9).Value = 1449  # ALC!I138
$ws.Cells.Item(138, 10).Value = 5800  # ALC!J138
$ws.Cells.Item(138, 11).Value = 4347  # ALC!K138
$ws.Cells.Item(138, 12).Value = 17400  # ALC!L138
$ws.Cells.Item(138, 13).Value = 793  # ALC!M138
$ws.Cells.Item(138, 14).Value = -27680  # ALC!N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1225.875  # ARM!H2
$ws.Cells.Item(2, 9).Value = 1225.875  # ARM!I2
$ws.Cells.Item(2, 11).Value = 1225.875  # ARM!K2
$ws.Cells.Item(2, 13).Value = -1112.875  # ARM!M2
$ws.Cells.Item(25, 8).Value = 2979  # ARM!H25
$ws.Cells.Item(25, 9).Value = 2979  # ARM!I25
$ws.Cells.Item(25, 11).Value = 2979  # ARM!K25
$ws.Cells.Item(25, 13).Value = -2577  # ARM!M25
$ws.Cells.Item(61, 8).Value = 2302.926  # ARM!H61
$ws.Cells.Item(61, 9).Value = 2009.579  # ARM!I61
$ws.Cells.Item(61, 10).Value = 2999.625  # ARM!J61
$ws.Cells.Item(61, 11).Value = 2009.579  # ARM!K61
$ws.Cells.Item(61, 12).Value = 2999.625  # ARM!L61
$ws.Cells.Item(61, 13).Value = -1797.579  # ARM!M61
$ws.Cells.Item(61, 14).Value = -3423.625  # ARM!N61
$ws.Cells.Item(74, 8).Value = 2170.9285  # ARM!H74
$ws.Cells.Item(74, 9).Value = 1167.5  # ARM!I74
$ws.Cells.Item(74, 10).Value = 2923.5  # ARM!J74
$ws.Cells.Item(74, 11).Value = 1167.5  # ARM!K74
$ws.Cells.Item(74, 12).Value = 2923.5  # ARM!L74
$ws.Cells.Item(74, 13).Value = -293.5  # ARM!M74
$ws.Cells.Item(74, 14).Value = -4671.5  # ARM!N74
$ws.Cells.Item(77, 8).Value = 2170.9285  # ARM!H77
$ws.Cells.Item(77, 9).Value = 1167.5  # ARM!I77
$ws.Cells.Item(77, 10).Value = 2923.5  # ARM!J77
$ws.Cells.Item(77, 11).Value = 5837.5  # ARM!K77
$ws.Cells.Item(77, 12).Value = 14617.5  # ARM!L77
$ws.Cells.Item(77, 13).Value = -1469.5  # ARM!M77
$ws.Cells.Item(77, 14).Value = -23353.5  # ARM!N77
$ws.Cells.Item(97, 8).Value = 849.2308  # ARM!H97
$ws.Cells.Item(97, 10).Value = 1241.3334  # ARM!J97
$ws.Cells.Item(97, 12).Value = 1241.3334  # ARM!L97
$ws.Cells.Item(97, 14).Value = -2233.3334  # ARM!N97
$ws.Cells.Item(102, 8).Value = 2316.7058  # ARM!H102
$ws.Cells.Item(102, 9).Value = 1573.75  # ARM!I102
$ws.Cells.Item(102, 11).Value = 1573.75  # ARM!K102
$ws.Cells.Item(102, 13).Value = 48.25  # ARM!M102
$ws.Cells.Item(116, 8).Value = 1225.875  # ARM!H116
$ws.Cells.Item(116, 9).Value = 1225.875  # ARM!I116
$ws.Cells.Item(116, 11).Value = 1225.875  # ARM!K116
$ws.Cells.Item(116, 13).Value = 1068.125  # ARM!M116
$ws.Cells.Item(122, 8).Value = 4001  # ARM!H122
$ws.Cells.Item(122, 9).Value = 4001  # ARM!I122
$ws.Cells.Item(122, 10).Value = 0  # ARM!J122
$ws.Cells.Item(122, 11).Value = 12003  # ARM!K122
$ws.Cells.Item(122, 12).Value = 0  # ARM!L122
$ws.Cells.Item(122, 13).Value = -9553  # ARM!M122
$ws.Cells.Item(122, 14).ClearContents()  # ARM!N122
$ws.Cells.Item(132, 8).Value = 2435.5  # ARM!H132
$ws.Cells.Item(132, 9).Value = 2063  # ARM!I132
$ws.Cells.Item(132, 10).Value = 2808  # ARM!J132
$ws.Cells.Item(132, 11).Value = 6189  # ARM!K132
$ws.Cells.Item(132, 12).Value = 8424  # ARM!L132
$ws.Cells.Item(132, 13).Value = -3659  # ARM!M132
$ws.Cells.Item(132, 14).Value = -13484  # ARM!N132
$ws.Cells.Item(136, 8).Value = 2302.926  # ARM!H136
$ws.Cells.Item(136, 9).Value = 2009.579  # ARM!I136
$ws.Cells.Item(136, 10).Value = 2999.625  # ARM!J136
$ws.Cells.Item(136, 11).Value = 6028.737  # ARM!K136
$ws.Cells.Item(136, 12).Value = 8998.875  # ARM!L136
$ws.Cells.Item(136, 13).Value = -3478.737  # ARM!M136
$ws.Cells.Item(136, 14).Value = -14098.875  # ARM!N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1225.875  # BSM!H3
$ws.Cells.Item(3, 9).Value = 1225.875  # BSM!I3
$ws.Cells.Item(3, 11).Value = 1225.875  # BSM!K3
$ws.Cells.Item(3, 13).Value = -1111.875  # BSM!M3
$ws.Cells.Item(10, 8).Value = 2000  # BSM!H10
$ws.Cells.Item(10, 9).Value = 2000  # BSM!I10
$ws.Cells.Item(10, 11).Value = 2000  # BSM!K10
$ws.Cells.Item(10, 13).Value = -1860  # BSM!M10
$ws.Cells.Item(12, 8).Value = 245  # BSM!H12
$ws.Cells.Item(12, 9).Value = 0  # BSM!I12
$ws.Cells.Item(12, 10).Value = 245  # BSM!J12
$ws.Cells.Item(12, 11).Value = 0  # BSM!K12
$ws.Cells.Item(12, 12).Value = 245  # BSM!L12
$ws.Cells.Item(12, 13).ClearContents()  # BSM!M12
$ws.Cells.Item(12, 14).Value = -581  # BSM!N12
$ws.Cells.Item(20, 8).Value = 71439130  # BSM!H20
$ws.Cells.Item(20, 9).Value = 83345500  # BSM!I20
$ws.Cells.Item(20, 10).Value = 899  # BSM!J20
$ws.Cells.Item(20, 11).Value = 83345500  # BSM!K20
$ws.Cells.Item(20, 12).Value = 899  # BSM!L20
$ws.Cells.Item(20, 13).Value = -83345253  # BSM!M20
$ws.Cells.Item(20, 14).Value = -1393  # BSM!N20
$ws.Cells.Item(25, 8).Value = 338.33334  # BSM!H25
$ws.Cells.Item(25, 9).Value = 338.33334  # BSM!I25
$ws.Cells.Item(25, 11).Value = 338.33334  # BSM!K25
$ws.Cells.Item(25, 13).Value = -103.33334  # BSM!M25
$ws.Cells.Item(134, 8).Value = 1476.1  # BSM!H134
$ws.Cells.Item(134, 9).Value = 863.7241  # BSM!I134
$ws.Cells.Item(134, 10).Value = 3090.5454  # BSM!J134
$ws.Cells.Item(134, 11).Value = 2591.1723  # BSM!K134
$ws.Cells.Item(134, 12).Value = 9271.636200000001  # BSM!L134
$ws.Cells.Item(134, 13).Value = -56.17230000000018  # BSM!M134
$ws.Cells.Item(134, 14).Value = -14341.6362  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 25005962  # CRP!H31
$ws.Cells.Item(31, 9).Value = 0  # CRP!I31
$ws.Cells.Item(31, 11).Value = 0  # CRP!K31
$ws.Cells.Item(31, 13).ClearContents()  # CRP!M31
$ws.Cells.Item(34, 8).Value = 25005962  # CRP!H34
$ws.Cells.Item(34, 9).Value = 0  # CRP!I34
$ws.Cells.Item(34, 11).Value = 0  # CRP!K34
$ws.Cells.Item(34, 13).ClearContents()  # CRP!M34
$ws.Cells.Item(58, 8).Value = 1374.6522  # CRP!H58
$ws.Cells.Item(58, 9).Value = 836.64703  # CRP!I58
$ws.Cells.Item(58, 11).Value = 836.64703  # CRP!K58
$ws.Cells.Item(58, 13).Value = -633.64703  # CRP!M58
$ws.Cells.Item(109, 8).Value = 49999  # CRP!H109
$ws.Cells.Item(109, 10).Value = 49999  # CRP!J109
$ws.Cells.Item(109, 12).Value = 49999  # CRP!L109
$ws.Cells.Item(109, 14).Value = -52079  # CRP!N109
$ws.Cells.Item(122, 8).Value = 3284.8235  # CRP!H122
$ws.Cells.Item(122, 9).Value = 2479.4  # CRP!I122
$ws.Cells.Item(122, 10).Value = 4435.4287  # CRP!J122
$ws.Cells.Item(122, 11).Value = 7438.200000000001  # CRP!K122
$ws.Cells.Item(122, 12).Value = 13306.2861  # CRP!L122
$ws.Cells.Item(122, 13).Value = -4988.200000000001  # CRP!M122
$ws.Cells.Item(122, 14).Value = -18206.2861  # CRP!N122
$ws.Cells.Item(132, 8).Value = 3999.0857  # CRP!H132
$ws.Cells.Item(132, 9).Value = 3107.76  # CRP!I132
$ws.Cells.Item(132, 11).Value = 9323.280000000001  # CRP!K132
$ws.Cells.Item(132, 13).Value = -6793.280000000001  # CRP!M132
$ws.Cells.Item(136, 8).Value = 1374.6522  # CRP!H136
$ws.Cells.Item(136, 9).Value = 836.64703  # CRP!I136
$ws.Cells.Item(136, 11).Value = 2509.94109  # CRP!K136
$ws.Cells.Item(136, 13).Value = 40.0589100000002  # CRP!M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 4999  # CUL!H80
$ws.Cells.Item(80, 9).Value = 0  # CUL!I80
$ws.Cells.Item(80, 10).Value = 4999  # CUL!J80
$ws.Cells.Item(80, 11).Value = 0  # CUL!K80
$ws.Cells.Item(80, 12).Value = 14997  # CUL!L80
$ws.Cells.Item(80, 13).ClearContents()  # CUL!M80
$ws.Cells.Item(80, 14).Value = -16869  # CUL!N80
$ws.Cells.Item(83, 8).Value = 4999  # CUL!H83
$ws.Cells.Item(83, 9).Value = 0  # CUL!I83
$ws.Cells.Item(83, 10).Value = 4999  # CUL!J83
$ws.Cells.Item(83, 11).Value = 0  # CUL!K83
$ws.Cells.Item(83, 12).Value = 44991  # CUL!L83
$ws.Cells.Item(83, 13).ClearContents()  # CUL!M83
$ws.Cells.Item(83, 14).Value = -54351  # CUL!N83

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 111113780  # GSM!H80
$ws.Cells.Item(80, 9).Value = 250002370  # GSM!I80
$ws.Cells.Item(80, 10).Value = 2899.6  # GSM!J80
$ws.Cells.Item(80, 11).Value = 250002370  # GSM!K80
$ws.Cells.Item(80, 12).Value = 2899.6  # GSM!L80
$ws.Cells.Item(80, 13).Value = -250001372  # GSM!M80
$ws.Cells.Item(80, 14).Value = -4895.6  # GSM!N80
$ws.Cells.Item(83, 8).Value = 111113780  # GSM!H83
$ws.Cells.Item(83, 9).Value = 250002370  # GSM!I83
$ws.Cells.Item(83, 10).Value = 2899.6  # GSM!J83
$ws.Cells.Item(83, 11).Value = 1250011850  # GSM!K83
$ws.Cells.Item(83, 12).Value = 14498  # GSM!L83
$ws.Cells.Item(83, 13).Value = -1250006858  # GSM!M83
$ws.Cells.Item(83, 14).Value = -24482  # GSM!N83
$ws.Cells.Item(92, 8).Value = 10000  # GSM!H92
$ws.Cells.Item(92, 10).Value = 10000  # GSM!J92
$ws.Cells.Item(92, 12).Value = 10000  # GSM!L92
$ws.Cells.Item(92, 14).Value = -13744  # GSM!N92
$ws.Cells.Item(132, 8).Value = 2036.9706  # GSM!H132
$ws.Cells.Item(132, 9).Value = 1898.909  # GSM!I132
$ws.Cells.Item(132, 10).Value = 2290.0833  # GSM!J132
$ws.Cells.Item(132, 11).Value = 5696.727000000001  # GSM!K132
$ws.Cells.Item(132, 12).Value = 6870.249899999999  # GSM!L132
$ws.Cells.Item(132, 13).Value = -3166.727000000001  # GSM!M132
$ws.Cells.Item(132, 14).Value = -11930.2499  # GSM!N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(23, 8).Value = 25013000  # LTW!H23
$ws.Cells.Item(23, 9).Value = 50000000  # LTW!I23
$ws.Cells.Item(23, 11).Value = 50000000  # LTW!K23
$ws.Cells.Item(23, 13).Value = -49999770  # LTW!M23
$ws.Cells.Item(122, 8).Value = 4412.5713  # LTW!H122
$ws.Cells.Item(122, 9).Value = 2722.5  # LTW!I122
$ws.Cells.Item(122, 10).Value = 6666  # LTW!J122
$ws.Cells.Item(122, 11).Value = 8167.5  # LTW!K122
$ws.Cells.Item(122, 12).Value = 19998  # LTW!L122
$ws.Cells.Item(122, 13).Value = -5717.5  # LTW!M122
$ws.Cells.Item(122, 14).Value = -24898  # LTW!N122
$ws.Cells.Item(123, 8).Value = 90428  # LTW!H123
$ws.Cells.Item(123, 10).Value = 90428  # LTW!J123
$ws.Cells.Item(123, 12).Value = 90428  # LTW!L123
$ws.Cells.Item(123, 14).Value = -100228  # LTW!N123

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 18511.834  # WVR!H41
$ws.Cells.Item(41, 9).Value = 17500  # WVR!I41
$ws.Cells.Item(41, 10).Value = 18714.2  # WVR!J41
$ws.Cells.Item(41, 11).Value = 17500  # WVR!K41
$ws.Cells.Item(41, 12).Value = 18714.2  # WVR!L41
$ws.Cells.Item(41, 13).Value = -17110  # WVR!M41
$ws.Cells.Item(41, 14).Value = -19494.2  # WVR!N41
$ws.Cells.Item(62, 8).Value = 7249.3335  # WVR!H62
$ws.Cells.Item(62, 9).Value = 1499  # WVR!I62
$ws.Cells.Item(62, 10).Value = 10124.5  # WVR!J62
$ws.Cells.Item(62, 11).Value = 1499  # WVR!K62
$ws.Cells.Item(62, 12).Value = 10124.5  # WVR!L62
$ws.Cells.Item(62, 13).Value = -875  # WVR!M62
$ws.Cells.Item(62, 14).Value = -11372.5  # WVR!N62
$ws.Cells.Item(65, 8).Value = 7249.3335  # WVR!H65
$ws.Cells.Item(65, 9).Value = 1499  # WVR!I65
$ws.Cells.Item(65, 10).Value = 10124.5  # WVR!J65
$ws.Cells.Item(65, 11).Value = 7495  # WVR!K65
$ws.Cells.Item(65, 12).Value = 50622.5  # WVR!L65
$ws.Cells.Item(65, 13).Value = -4375  # WVR!M65
$ws.Cells.Item(65, 14).Value = -56862.5  # WVR!N65
$ws.Cells.Item(96, 8).Value = 1750  # WVR!H96
$ws.Cells.Item(122, 8).Value = 15627330  # WVR!H122
$ws.Cells.Item(122, 9).Value = 2380.6365  # WVR!I122
$ws.Cells.Item(122, 11).Value = 7141.9095  # WVR!K122
$ws.Cells.Item(122, 13).Value = -4691.9095  # WVR!M122
